$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values in rows 2-5 (A:C) to reflect the randomized data
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = 2
$ws.Cells.Item(2, 3).Value = 4

$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = 4

$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = 4
$ws.Cells.Item(4, 3).Value = 5

$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 5

# Remove rows 6-7, which are no longer part of the dataset
$ws.Range("A6:C7").Delete()
